# The deck ships with two DrawingML theme parts:
#   theme1.xml -> "Office Theme" (the stock Office color palette) - currently
#                 wired only to the Notes Master, which this host doesn't
#                 expose for editing.
#   theme2.xml -> "Integral" / "Red Violet" palette - wired to the one and
#                 only Slide Master, i.e. the theme that actually drives
#                 every slide's look.
#
# The commit swaps the two palettes so the deck's live look becomes the
# plain "Office" palette. We reproduce that by rewriting every one of the
# 12 theme colors on the live (slide-master-backed) theme to the stock
# Office values, using the 12-slot ThemeColorScheme (maps 1:1 onto
# dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink in that order).

$p = $ppt.ActivePresentation

# RGB() isn't available in this host, so build COM color ints (0xBBGGRR)
# from hex strings ourselves.
function HexToCom([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette: the stock "Office" theme color scheme (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink), in that fixed order.
$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $themeColors.Item($i).RGB = HexToCom $officeColors[$i - 1]
}
